# ---------------------------------------------------------------------
# assignment05_key.docx edits:
#   1) Set NoProofing (w:noProof) on the two inline pictures near the top
#      of the document (question 4 screenshots / chart).
#   2) Rewrite the Q5 write-up paragraph to report the Pearson correlation
#      results instead of the old OLS-regression language.
# ---------------------------------------------------------------------

$d = $word.ActiveDocument

# --- 1) Mark the two figures as NoProof (adds <w:noProof/> to their run) ---
$d.InlineShapes.Item(1).Range.NoProofing = $true
$d.InlineShapes.Item(2).Range.NoProofing = $true

# --- 2) Replace the old OLS paragraph text with the new correlation text ---
$old = "I further fit an OLS regression model to estimate whether the increase in GPA is significantly related to the earliest age of the individual when she first felt being influenced by the mentor. The regression coefficient of age is 0.01 with a standard error of 0.01 ("
$rng = $d.Content
$rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$start = $rng.Start
$rng.Text = "My null hypothesis is there was no relationship between the mentee’s age and high school overall GPA. I conducted a Pearson’s correlation test in R to test this hypothesis. The correlation coefficient was 0.04 ("

# Re-apply (no-op) character formatting over each new sentence/phrase so Word
# keeps them as discrete runs, matching how the author authored the edit.
$r = $d.Range($start + 0, $start + 171)
$r.Bold = $true
$r.Bold = $false
$r = $d.Range($start + 171, $start + 177)
$r.Bold = $true
$r.Bold = $false
$r = $d.Range($start + 177, $start + 188)
$r.Bold = $true
$r.Bold = $false
$r = $d.Range($start + 188, $start + 201)
$r.Bold = $true
$r.Bold = $false
$r = $d.Range($start + 201, $start + 210)
$r.Bold = $true
$r.Bold = $false
$r = $d.Range($start + 210, $start + 211)
$r.Bold = $true
$r.Bold = $false

Write-Output "Applied noProof + Q5 rewrite edits"
